$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Title paragraph: "Θέμα: « Έγκριση ..." -> "Θέμα: «Έγκριση ..."
#    (drop the stray space right after the opening guillemet)
# ------------------------------------------------------------------
$d.Content.Find.Execute('Θέμα: « Έγκριση', $true, $false, $false, $false, $false, $true, 1, $false, 'Θέμα: «Έγκριση', 2) | Out-Null

# ------------------------------------------------------------------
# 2) Title paragraph: the country is no longer wrapped with a
#    hard-coded "στη χώρα «...»" phrase - the article/preposition now
#    comes from the database value itself, so just drop "στη χώρα «"
#    before ${country} and keep the placeholder unquoted.
#    The __DdeLink bookmark that used to start right before "χώρα"
#    must now start right after "school" (i.e. right before the
#    run that now reads "}» ") while still ending right after the
#    closing "}" of ${country}.
# ------------------------------------------------------------------

# Locate the end of "school" - that's where the bookmark must now start.
$schoolRng = $d.Content
$schoolRng.Find.Execute('school', $true, $false, $false, $false, $false, $true, 1, $false, '', 0) | Out-Null
$schoolEnd = $schoolRng.End

# Drop the old bookmark - it will be re-created in the right spot below.
$oldBookmarkName = '__DdeLink__157_983629326'
if ($d.Bookmarks.Exists($oldBookmarkName)) {
    $d.Bookmarks.Item($oldBookmarkName).Delete()
}

# Remove "στη χώρα «" between "}» " and "${country}".
$titleRng = $d.Content
$titleRng.Find.Execute('}» στη χώρα «${', $true, $false, $false, $false, $false, $true, 1, $false, '}» ${', 2) | Out-Null

# Re-find the closing brace of ${country} to know where the bookmark ends.
$countryRng = $d.Content
$countryRng.Find.Execute('${country}', $true, $false, $false, $false, $false, $true, 1, $false, '', 0) | Out-Null
$countryEnd = $countryRng.End

$bmRange = $d.Range($schoolEnd, $countryEnd)
$d.Bookmarks.Add($oldBookmarkName, $bmRange) | Out-Null

# ------------------------------------------------------------------
# 3) Second occurrence ("Εγκρίνουμε τη μετακίνηση ... προκειμένου να
#    μεταβούν στη χώρα «${country}» στο πλαίσιο ..."): drop the same
#    hard-coded "στη χώρα «" prefix and the trailing "»" so only the
#    bare ${country} placeholder remains.
# ------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute('μεταβούν στη χώρα «', $true, $false, $false, $false, $false, $true, 1, $false, 'μεταβούν ', 2) | Out-Null

$rng3 = $d.Content
$rng3.Find.Execute('country}»', $true, $false, $false, $false, $false, $true, 1, $false, 'country}', 2) | Out-Null
